# Refresh cryptocurrency Price / Volume(1h) figures (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.643.77"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "1.828.27"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "'309.29"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").Value = "'0.4682"
$ws.Range("E7").Value = "  +3.64%  "
$ws.Range("D8").Value = "'0.3597"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("D10").Value = "'0.9040"
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("D11").Value = "'0.07674"
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").Value = "'19.44"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "1.874.65"
$ws.Range("E13").Value = "  +4.15%  "
$ws.Range("D14").Value = "'5.265"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "'6.370"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "'87.65"
$ws.Range("E16").Value = "  +3.21%  "
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "'0.000008564"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "26.644.44"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").Value = "'14.22"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "'5.025"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").Value = "'10.55"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'1.910"
$ws.Range("E24").Value = "  -3.13%  "
$ws.Range("D25").Value = "'152.98"
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("D26").Value = "'17.93"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").Value = "'2.004"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").Value = "'113.69"
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("D31").Value = "'3.151"
$ws.Range("E31").Value = "  +2.02%  "
$ws.Range("D32").Value = "'2.849"
$ws.Range("E32").Value = "  +3.19%  "
$ws.Range("D33").Value = "'1.164"
$ws.Range("E33").Value = "  +5.62%  "
$ws.Range("D34").Value = "'0.7363"
$ws.Range("E34").Value = "  +1.93%  "
$ws.Range("D35").Value = "'4.431"
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("D37").Value = "'0.01933"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "'2.948"
$ws.Range("E38").Value = "  +3.29%  "
$ws.Range("D39").Value = "'0.05154"
$ws.Range("E39").Value = "  +1.34%  "
$ws.Range("D40").Value = "'6.869"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  -1.08%  "
$ws.Range("D43").Value = "'8.065"
$ws.Range("E43").Value = "  +0.60%  "
$ws.Range("D44").Value = "'1.009"
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").Value = "'0.4656"
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("D46").Value = "'10.06"
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("D47").Value = "'98.52"
$ws.Range("E47").Value = "  -2.14%  "
$ws.Range("D48").Value = "'1.573"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").Value = "'0.06022"
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("D50").Value = "'63.93"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").Value = "'35.84"
$ws.Range("E51").Value = "  -0.65%  "
